# Weekly update: insert two new daily price rows for "Femacal de La Calera - Papa"
# right above the existing block that starts at row 793. Excel's native Insert()
# shifts the old rows 793..873 down to 795..875, which is exactly what the
# published diff shows (every pre-existing row's content reappears two rows
# lower, unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the block (rows 793 and 794).
$ws.Range("A793:A794").EntireRow.Insert()

# --- New row 793 ---------------------------------------------------------
$ws.Cells.Item(793, 1).Value = 3
$ws.Cells.Item(793, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(793, 3).Value = "Coquimbo"
$ws.Cells.Item(793, 4).Value = 44946
$ws.Cells.Item(793, 5).Value = 5
$ws.Cells.Item(793, 6).Value = 100114001
$ws.Cells.Item(793, 7).Value = "Papa"
$ws.Cells.Item(793, 8).Value = "Rosara"
$ws.Cells.Item(793, 9).Value = "1a nueva(o)"
$ws.Cells.Item(793, 10).Value = 730
$ws.Cells.Item(793, 11).Value = 11000
$ws.Cells.Item(793, 12).Value = 12000
$ws.Cells.Item(793, 13).Value = 11479
$ws.Cells.Item(793, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(793, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(793, 16).Value = 459
$ws.Cells.Item(793, 17).Value = 25
$ws.Cells.Item(793, 18).Value = "Hortaliza"

# --- New row 794 ---------------------------------------------------------
$ws.Cells.Item(794, 1).Value = 3
$ws.Cells.Item(794, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(794, 3).Value = "Coquimbo"
$ws.Cells.Item(794, 4).Value = 44946
$ws.Cells.Item(794, 5).Value = 5
$ws.Cells.Item(794, 6).Value = 100114001
$ws.Cells.Item(794, 7).Value = "Papa"
$ws.Cells.Item(794, 8).Value = "Rosara"
$ws.Cells.Item(794, 9).Value = "2a nueva(o)"
$ws.Cells.Item(794, 10).Value = 120
$ws.Cells.Item(794, 11).Value = 10000
$ws.Cells.Item(794, 12).Value = 10000
$ws.Cells.Item(794, 13).Value = 10000
$ws.Cells.Item(794, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(794, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(794, 16).Value = 400
$ws.Cells.Item(794, 17).Value = 25
$ws.Cells.Item(794, 18).Value = "Hortaliza"

# Apply the same date style (numFmt) used by the rest of column D to the two
# new date cells, matching the s="2" attribute seen on every other row.
$ws.Cells.Item(793, 4).NumberFormat = $ws.Cells.Item(795, 4).NumberFormat
$ws.Cells.Item(794, 4).NumberFormat = $ws.Cells.Item(795, 4).NumberFormat
